$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (Fer Welcome Page task): Temps "1 h" -> "2h"
$ws.Range("D15").Value = "2h"

# Row 16 (Canviar spritesheets escenaris Handout task): fill in Persona Encarregada, Previsio, Temps
$ws.Range("B16").Value = "Rafa"
$ws.Range("C16").Value = "0:30 h"
$ws.Range("D16").Value = "1h"

# Row 18 (Fer End Battle Screen task): fill in Persona Encarregada, Previsio, Temps
$ws.Range("B18").Value = "Rafa"
$ws.Range("C18").Value = "0:30 h"
$ws.Range("D18").Value = "0:20 h"

# Update active selection to match final cursor position in the workbook
$ws.Range("D19").Select()
